$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.914.12"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").Value = "1.903.29"
$ws.Range("E3").Value = "  +0.80%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.8021"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "240.28"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.57%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3106"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "26.27"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07031"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07996"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.75%  "
$ws.Range("D12").Value = "1.913.39"
$ws.Range("E12").Value = "  +1.22%  "
$ws.Range("E13").Value = "  -0.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.168"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.35%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.19"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.01%  "
$ws.Range("D16").Value = "29.911.14"
$ws.Range("E16").Value = "  +0.47%  "
$ws.Range("E17").Value = "  +0.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.844"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.67"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("E20").Value = "  +1.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.003"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.18%  "
$ws.Range("D22").Value = "2.151.84"
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.898"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "168.53"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.174"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1464"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +15.54%  "
$ws.Range("E28").Value = "  +1.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.060"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.67%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.356"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.97%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.509"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.275"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05513"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.048"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.45%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.256"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7267"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.719"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01913"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.783"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.37%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4392"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "71.92"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.50%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.952"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.88%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.000"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8344"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.24%  "
$ws.Range("E45").Value = "  +0.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "100.78"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.539"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.689"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.76%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "2.059.84"
$ws.Range("E49").Value = "  +0.53%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "976.27"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.92%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "36.08"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
